# Updated cryptos list on Fri Jan 12 05:34:48 UTC 2024 with GitHub Actions
# Applies the price/volume refresh described by the diff.
#
# Column D ("Price") values must stay literal text (e.g. "46.110.44",
# "4.00", "108.70") rather than be auto-coerced by Excel into numbers
# (which would lose trailing zeros / use dotted-thousands notation /
# flip to scientific notation). Forcing the cell's NumberFormat to "@"
# (Text) before assigning the value keeps the exact text intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# Row 2 - Bitcoin
Set-PriceText "D2" "46.110.44"
$ws.Range("E2").Value = "  -0.45%  "

# Row 3 - Ethereum
Set-PriceText "D3" "2.594.87"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-PriceText "D5" "312.04"
$ws.Range("E5").Value = "  +1.60%  "

# Row 6 - Solana
Set-PriceText "D6" "99.11"
$ws.Range("E6").Value = "  -1.18%  "

# Row 7 - XRP
Set-PriceText "D7" "0.599"
$ws.Range("E7").Value = "  -0.32%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +2.08%  "

# Row 10 - Avalanche
Set-PriceText "D10" "39.09"
$ws.Range("E10").Value = "  +0.76%  "

# Row 11 - OKB
Set-PriceText "D11" "54.22"
$ws.Range("E11").Value = "  -1.73%  "

# Row 12 - Dogecoin
Set-PriceText "D12" "0.0842"
$ws.Range("E12").Value = "  +0.70%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  +0.06%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-PriceText "D14" "2.992.64"
$ws.Range("E14").Value = "  -0.34%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +1.61%  "

# Row 16 - WrappedEther
Set-PriceText "D16" "2.594.19"
$ws.Range("E16").Value = "  -1.28%  "

# Row 17 - Polygon
Set-PriceText "D17" "0.917"
$ws.Range("E17").Value = "  +2.07%  "

# Row 18 - Chainlink
Set-PriceText "D18" "14.86"
$ws.Range("E18").Value = "  -0.06%  "

# Row 19 - WrappedBTC
Set-PriceText "D19" "46.213.53"
$ws.Range("E19").Value = "  -0.55%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +1.31%  "

# Row 21 - Uniswap
Set-PriceText "D21" "6.75"
$ws.Range("E21").Value = "  +1.16%  "

# Row 22 - InternetComputer(DFINITY)
Set-PriceText "D22" "12.81"
$ws.Range("E22").Value = "  -3.33%  "

# Row 23 - BitcoinCash
Set-PriceText "D23" "293.79"
$ws.Range("E23").Value = "  +15.20%  "

# Row 24 - Litecoin
Set-PriceText "D24" "73.02"
$ws.Range("E24").Value = "  +2.67%  "

# Row 25 - PancakeSwap
Set-PriceText "D25" "3.07"
$ws.Range("E25").Value = "  +2.04%  "

# Row 26 - ImmutableX
$ws.Range("E26").Value = "  -0.26%  "

# Row 27 - EthereumClassic
Set-PriceText "D27" "29.49"
$ws.Range("E27").Value = "  +4.79%  "

# Row 28 - Dai
$ws.Range("E28").Value = "  +0.05%  "

# Row 29 - LEO
Set-PriceText "D29" "4.06"
$ws.Range("E29").Value = "  +1.26%  "

# Row 30 - Cosmos
Set-PriceText "D30" "10.84"
$ws.Range("E30").Value = "  +3.91%  "

# Row 31 - InjectiveProtocol
Set-PriceText "D31" "38.89"
$ws.Range("E31").Value = "  -3.00%  "

# Row 32 - Toncoin
$ws.Range("E32").Value = "  -2.65%  "

# Row 33 - Filecoin
Set-PriceText "D33" "6.22"
$ws.Range("E33").Value = "  +1.06%  "

# Row 34 - LidoDAOToken
Set-PriceText "D34" "3.59"
$ws.Range("E34").Value = "  -2.79%  "

# Row 35 - Monero
Set-PriceText "D35" "155.94"
$ws.Range("E35").Value = "  +2.61%  "

# Row 36 - Hedera
$ws.Range("E36").Value = "  +0.53%  "

# Row 37 - ARBITRUM
Set-PriceText "D37" "2.19"
$ws.Range("E37").Value = "  -6.15%  "

# Row 38 - WEMIXToken
Set-PriceText "D38" "2.78"
$ws.Range("E38").Value = "  -5.50%  "

# Row 39 - Kaspa
Set-PriceText "D39" "0.122"
$ws.Range("E39").Value = "  +4.03%  "

# Row 40 - Stellar
Set-PriceText "D40" "0.123"
$ws.Range("E40").Value = "  +1.41%  "

# Rows 41 & 42 swap places: Celestia moves to row 41, VeChain moves to row 42
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-PriceText "D41" "15.74"
$ws.Range("E41").Value = "  +0.85%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-PriceText "D42" "0.0332"
$ws.Range("E42").Value = "  +3.11%  "

# Row 43 - NEARProtocol
Set-PriceText "D43" "3.58"
$ws.Range("E43").Value = "  -0.25%  "

# Row 44 - RenderToken
Set-PriceText "D44" "4.00"
$ws.Range("E44").Value = "  -4.54%  "

# Row 45 - EnergySwap
$ws.Range("E45").Value = "  +8.56%  "

# Row 46 - Maker
Set-PriceText "D46" "2.107.22"
$ws.Range("E46").Value = "  +2.69%  "

# Row 47 - BitcoinSV
Set-PriceText "D47" "98.00"
$ws.Range("E47").Value = "  +7.59%  "

# Row 48 - FirstDigitalUSD
Set-PriceText "D48" "0.999"
$ws.Range("E48").Value = "  -0.02%  "

# Row 49 - FraxShare
Set-PriceText "D49" "9.49"
$ws.Range("E49").Value = "  +3.99%  "

# Row 50 - Aave
Set-PriceText "D50" "108.70"
$ws.Range("E50").Value = "  -0.69%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  +1.05%  "
